$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 5889.478
$ws.Range("I15").Value = 5889.478
$ws.Range("K15").Value = 17668.434
$ws.Range("M15").Value = -17499.434
$ws.Range("H17").Value = 2342.7834
$ws.Range("J17").Value = 2380.4575
$ws.Range("L17").Value = 7141.372499999999
$ws.Range("N17").Value = -7477.372499999999
$ws.Range("H19").Value = 330.69232
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 299.9091
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 299.9091
$ws.Range("M19").Value = -325
$ws.Range("N19").Value = -649.9091000000001
$ws.Range("H45").Value = 5863.4
$ws.Range("I45").Value = 2817
$ws.Range("J45").Value = 6625
$ws.Range("K45").Value = 8451
$ws.Range("L45").Value = 19875
$ws.Range("M45").Value = -8259
$ws.Range("N45").Value = -20259
$ws.Range("H112").Value = 2392.6875
$ws.Range("I112").Value = 587.5
$ws.Range("J112").Value = 3475.8
$ws.Range("K112").Value = 1762.5
$ws.Range("L112").Value = 10427.4
$ws.Range("M112").Value = -654.5
$ws.Range("N112").Value = -12643.4
$ws.Range("H129").Value = 1277.9773
$ws.Range("J129").Value = 1049.919
$ws.Range("L129").Value = 3149.757000000001
$ws.Range("N129").Value = -13149.757
$ws.Range("H137").Value = 3779.4878
$ws.Range("I137").Value = 960.7083
$ws.Range("J137").Value = 7758.9414
$ws.Range("K137").Value = 2882.1249
$ws.Range("L137").Value = 23276.8242
$ws.Range("M137").Value = -332.1248999999998
$ws.Range("N137").Value = -28376.8242
$ws.Range("H138").Value = 1498.12
$ws.Range("I138").Value = 701.02856
$ws.Range("J138").Value = 1927.3231
$ws.Range("K138").Value = 2103.08568
$ws.Range("L138").Value = 5781.969300000001
$ws.Range("M138").Value = 3036.91432
$ws.Range("N138").Value = -16061.9693
$ws.Range("H141").Value = 5156.722
$ws.Range("I141").Value = 1562.625
$ws.Range("J141").Value = 8032
$ws.Range("K141").Value = 4687.875
$ws.Range("L141").Value = 24096
$ws.Range("M141").Value = 492.125
$ws.Range("N141").Value = -34456

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2731.8
$ws.Range("I45").Value = 2644.8
$ws.Range("J45").Value = 2905.8
$ws.Range("K45").Value = 2644.8
$ws.Range("L45").Value = 2905.8
$ws.Range("M45").Value = -2267.8
$ws.Range("N45").Value = -3659.8
$ws.Range("H61").Value = 1293.091
$ws.Range("I61").Value = 944.7308
$ws.Range("K61").Value = 944.7308
$ws.Range("M61").Value = -732.7308
$ws.Range("H74").Value = 1880.2941
$ws.Range("I74").Value = 1586.8572
$ws.Range("K74").Value = 1586.8572
$ws.Range("M74").Value = -712.8571999999999
$ws.Range("H77").Value = 1880.2941
$ws.Range("I77").Value = 1586.8572
$ws.Range("K77").Value = 7934.286
$ws.Range("M77").Value = -3566.286
$ws.Range("H132").Value = 2110.2632
$ws.Range("I132").Value = 1183.3478
$ws.Range("J132").Value = 3531.5334
$ws.Range("K132").Value = 3550.0434
$ws.Range("L132").Value = 10594.6002
$ws.Range("M132").Value = -1020.0434
$ws.Range("N132").Value = -15654.6002
$ws.Range("H136").Value = 1293.091
$ws.Range("I136").Value = 944.7308
$ws.Range("K136").Value = 2834.1924
$ws.Range("M136").Value = -284.1923999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1454.82
$ws.Range("I107").Value = 1278.3235
$ws.Range("J107").Value = 1829.875
$ws.Range("K107").Value = 1278.3235
$ws.Range("L107").Value = 1829.875
$ws.Range("M107").Value = 641.6765
$ws.Range("N107").Value = -5669.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4310.1396
$ws.Range("I31").Value = 2033
$ws.Range("K31").Value = 2033
$ws.Range("M31").Value = -1738
$ws.Range("H34").Value = 4310.1396
$ws.Range("I34").Value = 2033
$ws.Range("K34").Value = 2033
$ws.Range("M34").Value = -1831
$ws.Range("H134").Value = 501754.53
$ws.Range("I134").Value = 1254.8334
$ws.Range("J134").Value = 3504752.8
$ws.Range("K134").Value = 3764.5002
$ws.Range("L134").Value = 10514258.4
$ws.Range("M134").Value = -1229.5002
$ws.Range("N134").Value = -10519328.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4345.4
$ws.Range("J70").Value = 4640.9414
$ws.Range("L70").Value = 13922.8242
$ws.Range("N70").Value = -14552.8242
$ws.Range("H73").Value = 4345.4
$ws.Range("J73").Value = 4640.9414
$ws.Range("L73").Value = 13922.8242
$ws.Range("N73").Value = -16106.8242
$ws.Range("H104").Value = 3999
$ws.Range("J104").Value = 3999
$ws.Range("L104").Value = 11997
$ws.Range("N104").Value = -17239
$ws.Range("H107").Value = 10982.737
$ws.Range("I107").Value = 10467.8
$ws.Range("J107").Value = 11554.889
$ws.Range("K107").Value = 31403.4
$ws.Range("L107").Value = 34664.667
$ws.Range("M107").Value = -29483.4
$ws.Range("N107").Value = -38504.667
$ws.Range("H113").Value = 3209.3618
$ws.Range("I113").Value = 4023.7585
$ws.Range("J113").Value = 1897.2778
$ws.Range("K113").Value = 12071.2755
$ws.Range("L113").Value = 5691.8334
$ws.Range("M113").Value = -9901.2755
$ws.Range("N113").Value = -10031.8334
$ws.Range("H121").Value = 104307.53
$ws.Range("I121").Value = 476.66666
$ws.Range("J121").Value = 115844.3
$ws.Range("K121").Value = 1429.99998
$ws.Range("L121").Value = 347532.9
$ws.Range("M121").Value = -119.9999800000001
$ws.Range("N121").Value = -350152.9
$ws.Range("H131").Value = 4261.7354
$ws.Range("J131").Value = 1506.0714
$ws.Range("L131").Value = 4518.2142
$ws.Range("N131").Value = -14598.2142
$ws.Range("H132").Value = 2526.0908
$ws.Range("I132").Value = 1385.7142
$ws.Range("J132").Value = 3058.2666
$ws.Range("K132").Value = 12471.4278
$ws.Range("L132").Value = 27524.3994
$ws.Range("M132").Value = -9941.427799999999
$ws.Range("N132").Value = -32584.3994
$ws.Range("H134").Value = 38188896
$ws.Range("I134").Value = 49098000
$ws.Range("J134").Value = 7033
$ws.Range("K134").Value = 147294000
$ws.Range("L134").Value = 21099
$ws.Range("M134").Value = -147288930
$ws.Range("N134").Value = -31239
$ws.Range("H140").Value = 602492
$ws.Range("I140").Value = 751990
$ws.Range("J140").Value = 4500
$ws.Range("K140").Value = 2255970
$ws.Range("L140").Value = 13500
$ws.Range("M140").Value = -2250790
$ws.Range("N140").Value = -23860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 18968.334
$ws.Range("J26").Value = 18968.334
$ws.Range("L26").Value = 18968.334
$ws.Range("N26").Value = -19528.334
$ws.Range("H50").Value = 18968.334
$ws.Range("J50").Value = 18968.334
$ws.Range("L50").Value = 18968.334
$ws.Range("N50").Value = -19964.334
$ws.Range("H53").Value = 46455
$ws.Range("J53").Value = 46455
$ws.Range("L53").Value = 46455
$ws.Range("N53").Value = -47717

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2348.7693
$ws.Range("I61").Value = 2321.2727
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 2321.2727
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -2119.2727
$ws.Range("N61").Value = -2904
$ws.Range("H82").Value = 2028.1818
$ws.Range("I82").Value = 1812.2222
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 1812.2222
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -1451.2222
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 2028.1818
$ws.Range("I85").Value = 1812.2222
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1812.2222
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -564.2221999999999
$ws.Range("N85").Value = -5496
$ws.Range("H100").Value = 2668.9
$ws.Range("I100").Value = 2654.3333
$ws.Range("K100").Value = 2654.3333
$ws.Range("M100").Value = -2113.3333
$ws.Range("H113").Value = 2348.7693
$ws.Range("I113").Value = 2321.2727
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 2321.2727
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -151.2727
$ws.Range("N113").Value = -6840
$ws.Range("H132").Value = 2152.481
$ws.Range("I132").Value = 1606.3455
$ws.Range("K132").Value = 4819.0365
$ws.Range("M132").Value = -2289.0365
